$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "In Translation"
$overview.Range("C3").Value = "In Translation"
$overview.Range("B4").Value = "In Translation"
$overview.Range("C4").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "In Translation"
$zhcn.Range("B4").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "In Translation"
$dede.Range("B4").Value = "In Translation"
